$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 0: ALC row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()

# Hunk 1: ALC row 98
$ws.Range("H98").Value = 901.0833
$ws.Range("I98").Value = 901.2727
$ws.Range("J98").Value = 899
$ws.Range("K98").Value = 901.2727
$ws.Range("L98").Value = 899
$ws.Range("M98").Value = 596.7273
$ws.Range("N98").Value = -3895

# Hunk 2: ALC row 112
$ws.Range("H112").Value = 2238.8
$ws.Range("I112").Value = 864.3333
$ws.Range("J112").Value = 2827.8572
$ws.Range("K112").Value = 2592.9999
$ws.Range("L112").Value = 8483.571599999999
$ws.Range("M112").Value = -1484.9999
$ws.Range("N112").Value = -10699.5716

# Hunk 3: ALC row 118
$ws.Range("H118").Value = 912.1429000000001
$ws.Range("I118").Value = 346.75
$ws.Range("J118").Value = 1666
$ws.Range("K118").Value = 1040.25
$ws.Range("L118").Value = 4998
$ws.Range("M118").Value = 616.75
$ws.Range("N118").Value = -8312

# Hunk 4: ALC row 122
$ws.Range("H122").Value = 901.0833
$ws.Range("I122").Value = 901.2727
$ws.Range("J122").Value = 899
$ws.Range("K122").Value = 2703.8181
$ws.Range("L122").Value = 2697
$ws.Range("M122").Value = -253.8181
$ws.Range("N122").Value = -7597

# Hunk 5: ALC row 127
$ws.Range("H127").Value = 2513.4
$ws.Range("I127").Value = 2589.3333
$ws.Range("J127").Value = 2399.5
$ws.Range("K127").Value = 7767.999899999999
$ws.Range("L127").Value = 7198.5
$ws.Range("M127").Value = -2807.999899999999
$ws.Range("N127").Value = -17118.5

# Hunk 6: ALC row 132
$ws.Range("H132").Value = 2210.818
$ws.Range("I132").Value = 2210.818
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6632.454000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4102.454000000001

# Hunk 7: ALC row 138
$ws.Range("H138").Value = 3051.1765
$ws.Range("I138").Value = 1412.5625
$ws.Range("J138").Value = 4507.722
$ws.Range("K138").Value = 4237.6875
$ws.Range("L138").Value = 13523.166
$ws.Range("M138").Value = 902.3125
$ws.Range("N138").Value = -23803.166

# Hunk 8: ALC row 141
$ws.Range("H141").Value = 2453.5
$ws.Range("I141").Value = 2236.182
$ws.Range("J141").Value = 4844
$ws.Range("K141").Value = 6708.545999999999
$ws.Range("L141").Value = 14532
$ws.Range("M141").Value = -1528.545999999999
$ws.Range("N141").Value = -24892

$ws = $wb.Worksheets.Item("ARM")
# Hunk 9: ARM row 132
$ws.Range("H132").Value = 1719.5834
$ws.Range("I132").Value = 1664
$ws.Range("J132").Value = 1997.5
$ws.Range("K132").Value = 4992
$ws.Range("L132").Value = 5992.5
$ws.Range("M132").Value = -2462
$ws.Range("N132").Value = -11052.5

$ws = $wb.Worksheets.Item("BSM")
# Hunk 10: BSM row 94
$ws.Range("H94").Value = 1344.5555
$ws.Range("I94").Value = 1593.9333
$ws.Range("J94").Value = 1032.8334
$ws.Range("K94").Value = 1593.9333
$ws.Range("L94").Value = 1032.8334
$ws.Range("M94").Value = -1142.9333
$ws.Range("N94").Value = -1934.8334

$ws = $wb.Worksheets.Item("CRP")
# Hunk 11: CRP row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()

# Hunk 12: CRP row 31
$ws.Range("H31").Value = 1051.3334
$ws.Range("I31").Value = 894.8
$ws.Range("J31").Value = 1247
$ws.Range("K31").Value = 894.8
$ws.Range("L31").Value = 1247
$ws.Range("M31").Value = -599.8
$ws.Range("N31").Value = -1837

# Hunk 13: CRP row 34
$ws.Range("H34").Value = 1051.3334
$ws.Range("I34").Value = 894.8
$ws.Range("J34").Value = 1247
$ws.Range("K34").Value = 894.8
$ws.Range("L34").Value = 1247
$ws.Range("M34").Value = -692.8
$ws.Range("N34").Value = -1651

# Hunk 14: CRP row 50
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

# Hunk 15: CRP row 55
$ws.Range("H55").Value = 37999.5
$ws.Range("I55").Value = 35000
$ws.Range("J55").Value = 38999.332
$ws.Range("K55").Value = 35000
$ws.Range("L55").Value = 38999.332
$ws.Range("M55").Value = -34685
$ws.Range("N55").Value = -39629.332

# Hunk 16: CRP row 86
$ws.Range("H86").Value = 3559.9167
$ws.Range("I86").Value = 3611
$ws.Range("J86").Value = 2998
$ws.Range("K86").Value = 3611
$ws.Range("L86").Value = 2998
$ws.Range("M86").Value = -2488
$ws.Range("N86").Value = -5244

# Hunk 17: CRP row 89
$ws.Range("H89").Value = 3559.9167
$ws.Range("I89").Value = 3611
$ws.Range("J89").Value = 2998
$ws.Range("K89").Value = 18055
$ws.Range("L89").Value = 14990
$ws.Range("M89").Value = -12439
$ws.Range("N89").Value = -26222

# Hunk 18: CRP row 132
$ws.Range("H132").Value = 4421.091
$ws.Range("I132").Value = 4421.091
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13263.273
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10733.273
$ws.Range("N132").ClearContents()

# Hunk 19: CRP row 134
$ws.Range("H134").Value = 1666.7778
$ws.Range("I134").Value = 1812.625
$ws.Range("J134").Value = 500
$ws.Range("K134").Value = 5437.875
$ws.Range("L134").Value = 1500
$ws.Range("M134").Value = -2902.875
$ws.Range("N134").Value = -6570

$ws = $wb.Worksheets.Item("CUL")
# Hunk 20: CUL row 80
$ws.Range("H80").Value = 6003
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 6003
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 18009
$ws.Range("N80").Value = -19881

# Hunk 21: CUL row 83
$ws.Range("H83").Value = 6003
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 6003
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 54027
$ws.Range("N83").Value = -63387

# Hunk 22: CUL row 93
$ws.Range("H93").Value = 19999
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 19999
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 59997
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -63741

# Hunk 23: CUL row 95
$ws.Range("H95").Value = 8865
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 8865
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 26595
$ws.Range("N95").Value = -30713

# Hunk 24: CUL row 123
$ws.Range("H123").Value = 5079.8
$ws.Range("I123").Value = 4133
$ws.Range("J123").Value = 6500
$ws.Range("K123").Value = 12399
$ws.Range("L123").Value = 19500
$ws.Range("M123").Value = -9949
$ws.Range("N123").Value = -24400

$ws = $wb.Worksheets.Item("GSM")
# Hunk 25: GSM row 20
$ws.Range("H20").Value = 33535000
$ws.Range("I20").Value = 50287500
$ws.Range("J20").Value = 29999
$ws.Range("K20").Value = 50287500
$ws.Range("L20").Value = 29999
$ws.Range("N20").Value = -30489
$ws.Range("M20").Value = -50287255

# Hunk 26: GSM row 80
$ws.Range("H80").Value = 3045
$ws.Range("I80").Value = 2787.2222
$ws.Range("J80").Value = 3625
$ws.Range("K80").Value = 2787.2222
$ws.Range("L80").Value = 3625
$ws.Range("M80").Value = -1789.2222
$ws.Range("N80").Value = -5621

# Hunk 27: GSM row 83
$ws.Range("H83").Value = 3045
$ws.Range("I83").Value = 2787.2222
$ws.Range("J83").Value = 3625
$ws.Range("K83").Value = 13936.111
$ws.Range("L83").Value = 18125
$ws.Range("M83").Value = -8944.111000000001
$ws.Range("N83").Value = -28109

# Hunk 28: GSM row 92
$ws.Range("H92").Value = 19643
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 19643
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 19643
$ws.Range("N92").Value = -23387

# Hunk 29: GSM row 102
$ws.Range("H102").Value = 2559.6
$ws.Range("I102").Value = 1699.5
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 1699.5
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -77.5
$ws.Range("N102").Value = -9244

$ws = $wb.Worksheets.Item("LTW")
# Hunk 30: LTW row 22
$ws.Range("H22").Value = 736.8125
$ws.Range("I22").Value = 670.9167
$ws.Range("J22").Value = 934.5
$ws.Range("K22").Value = 670.9167
$ws.Range("L22").Value = 934.5
$ws.Range("M22").Value = -375.9167
$ws.Range("N22").Value = -1524.5

# Hunk 31: LTW row 27
$ws.Range("H27").Value = 736.8125
$ws.Range("I27").Value = 670.9167
$ws.Range("J27").Value = 934.5
$ws.Range("K27").Value = 670.9167
$ws.Range("L27").Value = 934.5
$ws.Range("M27").Value = -563.9167
$ws.Range("N27").Value = -1148.5

# Hunk 32: LTW row 40
$ws.Range("H40").Value = 3043.2222
$ws.Range("I40").Value = 2150
$ws.Range("J40").Value = 3298.4285
$ws.Range("K40").Value = 2150
$ws.Range("L40").Value = 3298.4285
$ws.Range("M40").Value = -2014
$ws.Range("N40").Value = -3570.4285

# Hunk 33: LTW row 46
$ws.Range("H46").Value = 1645.7778
$ws.Range("I46").Value = 1441.6666
$ws.Range("J46").Value = 2666.3333
$ws.Range("K46").Value = 1441.6666
$ws.Range("L46").Value = 2666.3333
$ws.Range("M46").Value = -1253.6666
$ws.Range("N46").Value = -3042.3333

# Hunk 34: LTW row 64
$ws.Range("H64").Value = 22500
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 22500
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 22500
$ws.Range("N64").Value = -22950

# Hunk 35: LTW row 67
$ws.Range("H67").Value = 22500
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 22500
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 22500
$ws.Range("N67").Value = -24060

# Hunk 36: LTW row 82
$ws.Range("H82").Value = 1861
$ws.Range("I82").Value = 1872.5
$ws.Range("J82").Value = 1849.5
$ws.Range("K82").Value = 1872.5
$ws.Range("L82").Value = 1849.5
$ws.Range("M82").Value = -1511.5
$ws.Range("N82").Value = -2571.5

# Hunk 37: LTW row 85
$ws.Range("H85").Value = 1861
$ws.Range("I85").Value = 1872.5
$ws.Range("J85").Value = 1849.5
$ws.Range("K85").Value = 1872.5
$ws.Range("L85").Value = 1849.5
$ws.Range("M85").Value = -624.5
$ws.Range("N85").Value = -4345.5

# Hunk 38: LTW row 105
$ws.Range("H105").Value = 19998.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 19998.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 19998.5
$ws.Range("N105").Value = -26986.5

# Hunk 39: LTW row 106
$ws.Range("H106").Value = 20397.8
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 20397.8
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 20397.8
$ws.Range("N106").Value = -22921.8

# Hunk 40: LTW row 122
$ws.Range("H122").Value = 6521.839
$ws.Range("I122").Value = 4741.25
$ws.Range("J122").Value = 7646.421
$ws.Range("K122").Value = 14223.75
$ws.Range("L122").Value = 22939.263
$ws.Range("M122").Value = -11773.75
$ws.Range("N122").Value = -27839.263

# Hunk 41: LTW row 136
$ws.Range("H136").Value = 1617.7646
$ws.Range("I136").Value = 1409.3636
$ws.Range("J136").Value = 1999.8334
$ws.Range("K136").Value = 4228.0908
$ws.Range("L136").Value = 5999.5002
$ws.Range("M136").Value = -1678.0908
$ws.Range("N136").Value = -11099.5002

$ws = $wb.Worksheets.Item("WVR")
# Hunk 42: WVR row 15
$ws.Range("H15").Value = 38099
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 38099
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 38099
$ws.Range("N15").Value = -38675

# Hunk 43: WVR row 42
$ws.Range("H42").Value = 49999
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 49999
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 49999
$ws.Range("N42").Value = -50755

# Hunk 44: WVR row 63
$ws.Range("H63").Value = 32833
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 32833
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 32833
$ws.Range("N63").Value = -34081

# Hunk 45: WVR row 66
$ws.Range("H66").Value = 32833
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 32833
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 98499
$ws.Range("N66").Value = -104739

# Hunk 46: WVR row 68
$ws.Range("H68").Value = 50000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 50000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51622

# Hunk 47: WVR row 71
$ws.Range("H71").Value = 50000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 50000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -158112

# Hunk 48: WVR row 107
$ws.Range("H107").Value = 1799.8572
$ws.Range("I107").Value = 1520.2
$ws.Range("J107").Value = 2499
$ws.Range("K107").Value = 4560.6
$ws.Range("L107").Value = 7497
$ws.Range("M107").Value = -2640.6
$ws.Range("N107").Value = -11337

# Hunk 49: WVR row 122
$ws.Range("H122").Value = 7663.6
$ws.Range("I122").Value = 6323.5
$ws.Range("J122").Value = 8557
$ws.Range("K122").Value = 18970.5
$ws.Range("L122").Value = 25671
$ws.Range("M122").Value = -16520.5
$ws.Range("N122").Value = -30571

# Hunk 50: WVR row 132
$ws.Range("H132").Value = 1676.8572
$ws.Range("I132").Value = 1867.6
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 5602.799999999999
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -3072.799999999999
$ws.Range("N132").Value = -8660

# Hunk 51: WVR row 136
$ws.Range("H136").Value = 2234.7856
$ws.Range("I136").Value = 2507.9546
$ws.Range("J136").Value = 1233.1666
$ws.Range("K136").Value = 7523.8638
$ws.Range("L136").Value = 3699.4998
$ws.Range("M136").Value = -4973.8638
$ws.Range("N136").Value = -8799.4998
